# [FEATURE] Kyw Habilitar Command Line
# Add two new users to the "Users" sheet: CRECEREM / Usuario Emergencia
# and F02971 / Usuario sin command line. Widen column C to fit the new
# (longer) text and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Column C (description/role) already carries the right text style on the
# existing rows (right-aligned, text number format) - reuse it for the new
# rows instead of re-building it, so no new style gets added.
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C13:C14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New row 13: CRECEREM / Usuario Emergencia
$ws.Range("A13").Value = "CRECEREM"
$ws.Range("C13").Value = "Usuario Emergencia"

# New row 14: F02971 / Usuario sin command line
$ws.Range("A14").Value = "F02971"
$ws.Range("C14").Value = "Usuario sin command line"

# Widen column C so the longer descriptions fit.
$ws.Columns.Item(3).ColumnWidth = 23.7

# Move the active selection.
$ws.Range("E12").Select() | Out-Null
